$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# summer 24 week 13 inputs
$ws.Range("D2").Value = 1.27

$ws.Range("B4").Value = 1.46
$ws.Range("C4").Value = 1.45
$ws.Range("D4").Value = 1.33
$ws.Range("E4").Value = 1.22

$ws.Range("G6").Value = 0.97

$ws.Range("F7").Value = 1.47
